# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/centered/bordered) onto the new header
# cells, then set the header text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @{Row=2; I=1; J=7}
    @{Row=3; I=1; J=7}
    @{Row=4; I=1; J=6}
    @{Row=5; I=1; J=6}
    @{Row=6; I=5; J=7}
    @{Row=7; I=5; J=8}
    @{Row=8; I=2; J=6}
    @{Row=9; I=1; J=4}
    @{Row=10; I=5; J=8}
    @{Row=11; I=2; J=6}
    @{Row=12; I=6; J=7}
    @{Row=13; I=6; J=6}
    @{Row=14; I=9; J=9}
    @{Row=15; I=3; J=8}
    @{Row=16; I=2; J=7}
    @{Row=17; I=3; J=5}
    @{Row=18; I=5; J=6}
    @{Row=19; I=5; J=6}
    @{Row=20; I=4; J=6}
    @{Row=21; I=7; J=7}
    @{Row=22; I=5; J=5}
    @{Row=23; I=8; J=8}
    @{Row=24; I=3; J=6}
    @{Row=25; I=5; J=9}
    @{Row=26; I=3; J=7}
    @{Row=27; I=1; J=5}
    @{Row=28; I=5; J=6}
    @{Row=29; I=2; J=5}
    @{Row=30; I=4; J=5}
    @{Row=31; I=1; J=3}
    @{Row=32; I=1; J=5}
    @{Row=33; I=1; J=3}
    @{Row=34; I=2; J=3}
    @{Row=35; I=1; J=4}
    @{Row=36; I=1; J=4}
    @{Row=37; I=1; J=6}
    @{Row=38; I=7; J=8}
    @{Row=39; I=1; J=6}
    @{Row=40; I=5; J=8}
    @{Row=41; I=1; J=4}
    @{Row=42; I=1; J=3}
    @{Row=43; I=1; J=2}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
}
